$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" header in H1, reusing the same formatting (style) as
# the neighboring header cell G1 ("sum").
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)   # xlPasteFormats

# Add the value for the new "Save" column in the data row (H2), unstyled
# like the other data cells in that row.
$ws.Range("H2").Value = 1
